# Updated Data Access from Excel sheet for MSO Tests
# Append the new Sales-Order related lookup rows to the DataSheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSheet")

# New rows of (Key, Value) pairs appended below the existing A1:B4 block.
# Row 11 (SuspendSO / YE00194) is entered value-column-first so the shared
# string table grows in the original authoring order.
$ws.Cells.Item(5, 1).Value = "SOEnquire"
$ws.Cells.Item(5, 2).Value = "A00071A"

$ws.Cells.Item(6, 1).Value = "RouteSO"

$ws.Cells.Item(7, 1).Value = "RedirectSO"

$ws.Cells.Item(8, 1).Value = "BrowseSO"

$ws.Cells.Item(9, 1).Value = "AssignSO"
$ws.Cells.Item(9, 2).Value = "ZA00103"

$ws.Cells.Item(10, 1).Value = "SplitSO"
$ws.Cells.Item(10, 2).Value = "WI00204"

$ws.Cells.Item(11, 2).Value = "YE00194"
$ws.Cells.Item(11, 1).Value = "SuspendSO"

$ws.Cells.Item(12, 1).Value = "ManualWaitlistSO"

$ws.Cells.Item(13, 1).Value = "WaitlistSO"

$ws.Cells.Item(14, 1).Value = "RejectSO"

$ws.Cells.Item(15, 1).Value = "AccountSO"

$ws.Cells.Item(16, 1).Value = "GeneralSO"

# Move/collapse the selection to the last populated cell, matching the
# saved worksheet view state after the edit.
$ws.Range("A16").Select()
